$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Corrected rates for last 3 months (rows 112-114: Dec, Nov, Oct 2024)
$ws.Range("C112").Value = 4.33
$ws.Range("D112").Value = 4.11
$ws.Range("E112").Value = 3.86
$ws.Range("F112").Value = 3.76
$ws.Range("G112").Value = 3.51

$ws.Range("C113").Value = 4.32
$ws.Range("D113").Value = 3.87
$ws.Range("E113").Value = 3.59
$ws.Range("F113").Value = 3.54
$ws.Range("G113").Value = 3.17

$ws.Range("C114").Value = 4.08
$ws.Range("D114").Value = 3.74
$ws.Range("E114").Value = 3.31
$ws.Range("F114").Value = 3.28
$ws.Range("G114").Value = 3.17

# A new color-scale conditional-formatting rule got added covering one more
# row (C115:G115) below the existing per-row rules, bumping every existing
# rule's priority by one and inserting the new rule at priority 1 (top).
$newRule = $ws.Range("C115:G115").FormatConditions.AddColorScale(3)
$newRule.ColorScaleCriteria(1).FormatColor.Color = 8109667
$newRule.ColorScaleCriteria(2).FormatColor.Color = 8711167
$newRule.ColorScaleCriteria(3).FormatColor.Color = 7039480
$newRule.SetFirstPriority()

# Move the active selection to D119.
[void]$ws.Range("D119").Select()
